$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format before assignment so that
# numeric-looking strings (e.g. "245.17") are stored as literal text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.075.94'
$ws.Range("E2").Value = '  +4.88%  '
$ws.Range("D3").Value = '2.246.68'
$ws.Range("E3").Value = '  +4.07%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '245.17'
$ws.Range("E5").Value = '  +3.96%  '
$ws.Range("E6").Value = '  +1.91%  '
$ws.Range("D7").Value = '75.19'
$ws.Range("E7").Value = '  +8.85%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '0.605'
$ws.Range("E9").Value = '  +7.22%  '
$ws.Range("D10").Value = '40.98'
$ws.Range("E10").Value = '  +6.02%  '
$ws.Range("E11").Value = '  +2.62%  '
$ws.Range("D12").Value = '6.94'
$ws.Range("E12").Value = '  +4.47%  '
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("D14").Value = '2.586.71'
$ws.Range("E14").Value = '  +4.48%  '
$ws.Range("D15").Value = '14.60'
$ws.Range("E15").Value = '  +3.13%  '
$ws.Range("D16").Value = '2.241.70'
$ws.Range("E16").Value = '  +4.58%  '
$ws.Range("D17").Value = '0.794'
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("D18").Value = '42.980.91'
$ws.Range("E18").Value = '  +5.31%  '
$ws.Range("D19").Value = '0.0000104'
$ws.Range("E19").Value = '  +5.86%  '
$ws.Range("D20").Value = '71.14'
$ws.Range("E20").Value = '  +2.16%  '
$ws.Range("E21").Value = '  +3.78%  '
$ws.Range("D22").Value = '9.94'
$ws.Range("E22").Value = '  +7.18%  '
$ws.Range("D23").Value = '230.08'
$ws.Range("E23").Value = '  +2.63%  '
$ws.Range("D24").Value = '2.19'
$ws.Range("E24").Value = '  +16.67%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").Value = '10.87'
$ws.Range("E26").Value = '  +2.65%  '
$ws.Range("D27").Value = '3.42'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").Value = '2.24'
$ws.Range("E28").Value = '  +3.07%  '
$ws.Range("D29").Value = '38.72'
$ws.Range("E29").Value = '  +29.14%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '172.76'
$ws.Range("E30").Value = '  +2.94%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '2.13'
$ws.Range("E31").Value = '  -1.75%  '
$ws.Range("D32").Value = '20.30'
$ws.Range("E32").Value = '  +2.85%  '
$ws.Range("D33").Value = '0.0795'
$ws.Range("E33").Value = '  +6.31%  '
$ws.Range("D34").Value = '5.29'
$ws.Range("E34").Value = '  +4.26%  '
$ws.Range("E35").Value = '  +2.26%  '
$ws.Range("E36").Value = '  +7.86%  '
$ws.Range("D37").Value = '4.35'
$ws.Range("E37").Value = '  +7.39%  '
$ws.Range("D38").Value = '0.0332'
$ws.Range("E38").Value = '  +19.60%  '
$ws.Range("D39").Value = '13.09'
$ws.Range("E39").Value = '  +13.91%  '
$ws.Range("E40").Value = '  +4.15%  '
$ws.Range("D41").Value = '5.49'
$ws.Range("E41").Value = '  +3.98%  '
$ws.Range("E42").Value = '  +10.57%  '
$ws.Range("D43").Value = '59.54'
$ws.Range("E43").Value = '  +2.92%  '
$ws.Range("D44").Value = '105.23'
$ws.Range("E44").Value = '  +9.16%  '
$ws.Range("D45").Value = '8.70'
$ws.Range("E45").Value = '  +6.21%  '
$ws.Range("D46").Value = '0.479'
$ws.Range("E46").Value = '  +30.11%  '
$ws.Range("D47").Value = '0.0991'
$ws.Range("E47").Value = '  +3.90%  '
$ws.Range("D48").Value = '2.38'
$ws.Range("E48").Value = '  +10.82%  '
$ws.Range("E49").Value = '  +4.18%  '
$ws.Range("E50").Value = '  +3.78%  '
$ws.Range("D51").Value = '2.460.07'
$ws.Range("E51").Value = '  +4.35%  '

# Restore default (unstyled) formatting on column D so the cells
# match the original workbook formatting (no explicit style index).
$ws.Range("D2:D51").ClearFormats()
